$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 4 (Peso Muerto) Visible column: No -> Sí
$ws.Range("C4").Value = "Sí"

# New row 5 / row 6 content (entered column-A-first for both rows, then column B bottom-up)
$ws.Range("A5").Value = "100 metros "
$ws.Range("A6").Value = "400 metros"
$ws.Range("B6").Value = "400 Mts"
$ws.Range("B5").Value = "100 Mts"
$ws.Range("C5").Value = "Sí"
$ws.Range("C6").Value = "Sí"

# Row 7: Course Navette / Course Navette / Sí
$ws.Range("A7").Value = "Course Navette"
$ws.Range("B7").Value = "Course Navette"
$ws.Range("C7").Value = "Sí"

# Row 8: Dinamometria / Dinamometria Der. / Sí
$ws.Range("A8").Value = "Dinamometria"
$ws.Range("B8").Value = "Dinamometria Der."
$ws.Range("C8").Value = "Sí"

# Row 9: Dinamometria / Dinamometria Izq. / Sí
$ws.Range("A9").Value = "Dinamometria"
$ws.Range("B9").Value = "Dinamometria Izq."
$ws.Range("C9").Value = "Sí"

# Adjust column B width (auto-fit to contents, as Excel does with bestFit columns)
$ws.Columns.Item(2).ColumnWidth = 17.2

# Set the selection to C10 as in the saved file
$ws.Range("C10").Select()
